$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,('بازار ابزارهاي مشتقه', 'IFB.Derivatives')
    ,('بازار ابزارهاي مشتقه فرابورس', 'IFB.Derivatives')
    ,('بازار ابزارهای مشتقه', 'IFB.Derivatives')
    ,('بازار ابزارهای مشتقه فرابورس', 'IFB.Derivatives')
    ,('بازار اول فرابورس', 'IFB.M1')
    ,('بازار اول فرابورس لغو پذیرش شده', 'IFB.M1')
    ,('بازار دوم فرابورس', 'IFB.M2')
    ,('بازار دوم فرابورس قانون احکام دائمی برنامه های توسعه کشور', 'IFB.M2')
    ,('بازار دوم فرابورس قانون احکام دایمی برنامه های توسعه کشور', 'IFB.M2')
    ,('بازار دوم فرابورس لغو پذیرش شده', 'IFB.M2')
    ,('بازار سوم فرابورس', 'IFB.M3')
    ,('بازار سوم فرابورس قانون احکام دائمی برنامه های توسعه کشور', 'IFB.M3')
    ,('بازار سوم فرابورس قانون احکام دایمی برنامه های توسعه کشور', 'IFB.M3')
    ,('بازار ابزارهاي نوين مالي فرابورس', 'IFB.NFinTools')
    ,('بازار ابزارهای نوین مالی فرابورس', 'IFB.NFinTools')
    ,('بازار پايه  فرابورس', 'IFB.Paye')
    ,('بازار پايه  فرابورس قانون احکام دائمی برنامه های توسعه کشور', 'IFB.Paye')
    ,('بازار پایه فرابورس', 'IFB.Paye')
    ,('بازار پایه فرابورس قانون احکام دایمی برنامه های توسعه کشور', 'IFB.Paye')
    ,('بازار پايه نارنجي فرابورس', 'IFB.Paye.Orange')
    ,('بازار پايه نارنجي فرابورس قانون احکام دائمی برنامه های توسعه کشور', 'IFB.Paye.Orange')
    ,('بازار پايه نارنجي فرابورس قانون احکام دائمی برنامه های توسعه کشور/عدم اظهارنظر', 'IFB.Paye.Orange')
    ,('بازار پايه نارنجي فرابورس لغو پذیرش شده', 'IFB.Paye.Orange')
    ,('بازار پايه نارنجي فرابورس لغو پذیرش شده/ اظهارنظر مردود', 'IFB.Paye.Orange')
    ,('بازار پايه نارنجي فرابورس لغو پذیرش شده/ عدم اظهارنظر', 'IFB.Paye.Orange')
    ,('بازار پايه نارنجي فرابورس لغو پذیرش شده/عدم اظهارنظر', 'IFB.Paye.Orange')
    ,('بازار پایه نارنجی فرابورس', 'IFB.Paye.Orange')
    ,('بازار پایه نارنجی فرابورس قانون احکام دایمی برنامه های توسعه کشور', 'IFB.Paye.Orange')
    ,('بازار پایه نارنجی فرابورس قانون احکام دایمی برنامه های توسعه کشور/عدم اظهارنظر', 'IFB.Paye.Orange')
    ,('بازار پایه نارنجی فرابورس لغو پذیرش شده', 'IFB.Paye.Orange')
    ,('بازار پایه نارنجی فرابورس لغو پذیرش شده/ اظهارنظر مردود', 'IFB.Paye.Orange')
    ,('بازار پایه نارنجی فرابورس لغو پذیرش شده/ عدم اظهارنظر', 'IFB.Paye.Orange')
    ,('بازار پایه نارنجی فرابورس لغو پذیرش شده/عدم اظهارنظر', 'IFB.Paye.Orange')
    ,('بازار پايه قرمز فرابورس', 'IFB.Paye.Red')
    ,('بازار پايه قرمز فرابورس قانون احکام دائمی برنامه های توسعه کشور', 'IFB.Paye.Red')
    ,('بازار پايه قرمز فرابورس قانون احکام دائمی برنامه های توسعه کشور/ اظهارنظر مردود', 'IFB.Paye.Red')
    ,('بازار پايه قرمز فرابورس لغو پذیرش شده', 'IFB.Paye.Red')
    ,('بازار پايه قرمز فرابورس لغو پذیرش شده/ عدم اظهارنظر/ انحلال', 'IFB.Paye.Red')
    ,('بازار پايه قرمز فرابورس لغو پذیرش شده/اظهارنظر مردود', 'IFB.Paye.Red')
    ,('بازار پايه قرمز فرابورس لغو پذیرش شده/عدم اظهارنظر/ورشکستگی', 'IFB.Paye.Red')
    ,('بازار پایه قرمز فرابورس', 'IFB.Paye.Red')
    ,('بازار پایه قرمز فرابورس قانون احکام دایمی برنامه های توسعه کشور', 'IFB.Paye.Red')
    ,('بازار پایه قرمز فرابورس قانون احکام دایمی برنامه های توسعه کشور/ اظهارنظر مردود', 'IFB.Paye.Red')
    ,('بازار پایه قرمز فرابورس لغو پذیرش شده', 'IFB.Paye.Red')
    ,('بازار پایه قرمز فرابورس لغو پذیرش شده/ عدم اظهارنظر/ انحلال', 'IFB.Paye.Red')
    ,('بازار پایه قرمز فرابورس لغو پذیرش شده/اظهارنظر مردود', 'IFB.Paye.Red')
    ,('بازار پایه قرمز فرابورس لغو پذیرش شده/عدم اظهارنظر/ورشکستگی', 'IFB.Paye.Red')
    ,('بازار پايه زرد فرابورس', 'IFB.Paye.Yellow')
    ,('بازار پايه زرد فرابورس قانون احکام دائمی برنامه های توسعه کشور', 'IFB.Paye.Yellow')
    ,('بازار پايه زرد فرابورس لغو پذیرش شده', 'IFB.Paye.Yellow')
    ,('بازار پایه زرد فرابورس', 'IFB.Paye.Yellow')
    ,('بازار پایه زرد فرابورس قانون احکام دایمی برنامه های توسعه کشور', 'IFB.Paye.Yellow')
    ,('بازار پایه زرد فرابورس لغو پذیرش شده', 'IFB.Paye.Yellow')
    ,('شرکتهاي کوچک و متوسط فرابورس', 'IFB.SME')
    ,('شرکتهای کوچک و متوسط فرابورس', 'IFB.SME')
    ,('بورس کالا', 'IME')
    ,('بازار اوراق بدهي', 'TSE.Bonds')
    ,('بازار اوراق بدهی', 'TSE.Bonds')
    ,('بازار عادي آتي', 'TSE.Derivatives')
    ,('بازار عادی آتی', 'TSE.Derivatives')
    ,('بازار اول (تابلوي اصلي) بورس', 'TSE.M1.Main')
    ,('بازار اول (تابلوی اصلی) بورس', 'TSE.M1.Main')
    ,('بازار اول (تابلوي فرعي) بورس', 'TSE.M1.Subsidiary')
    ,('بازار اول (تابلوی فرعی) بورس', 'TSE.M1.Subsidiary')
    ,('بازار دوم بورس', 'TSE.M2')
    ,('بازار چهارم بورس', 'TSE.M4')
    ,('لغو پذیرش شده', '')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

Write-Output "Applied MarketTitle/MarketId reorder + text fixes to rows 2-68"